$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells that reuse already-existing shared strings (order among themselves
# does not affect the shared string table, since the strings already exist).
$ws.Range("D6").Value = "鐵捲門 一般門 窗戶"
$ws.Range("D14").Value = "可以刷卡"
$ws.Range("D17").Value = "印章(契約要圖面 要畫個圖"

# Cells introducing brand new shared strings. These must be written in the
# same order they were newly added to the shared string table so the
# resulting indices line up with the target file.
$ws.Range("D9").Value = "固定基地台(可支援sim卡"
$ws.Range("D15").Value = "一天"
$ws.Range("D16").Value = "有"
$ws.Range("D19").Value = "男的 洪啟煌 0933 185 241"
$ws.Range("D20").Value = "備註: 很愛聊天…電話費QQ"
$ws.Range("D18").Value = "有線無線都可以"
$ws.Range("D13").Value = "1200*11 / 1500*11"
$ws.Range("D12").Value = "就先聯絡畫圖"
$ws.Range("D10").Value = "有認識的 但可能沒比較便宜"

# Reflect the final selection left behind in the saved workbook.
$ws.Range("D24").Select()
